# Apply the "break out stock.yaml completed" update:
#  1. Append 5 new data rows (195-199) to the "day" sheet.
#  2. Fix the bsecode column (D) on the "week" sheet for rows 97-102 so it
#     is stored as a genuine number rather than a numeric-looking string.

$wb = $excel.ActiveWorkbook

# --- 1. "day" sheet: append new rows 195-199 ------------------------------
$dayWs = $wb.Worksheets.Item("day")

# Note: bsecode (column D) is written exactly as scraped from the source
# feed -- a numeric-looking piece of text, not a number -- so it is entered
# with a leading apostrophe to force Excel to keep it as text (matching the
# inlineStr cell type in the target workbook).
$newRows = @(
    @(1, "SIEMENS",  "Siemens Limited",             "'500550", 2.52,  7023.15, 373486,    "day", "22/07/2024 11:36:32"),
    @(2, "DIVISLAB", "Divi's Laboratories Limited", "'532488", -0.02, 4519.05, 308100,    "day", "22/07/2024 11:36:32"),
    @(3, "AARTIIND", "Aarti Industries Limited",    "'524208", 3.92,  697.65,  1018933,   "day", "22/07/2024 11:36:32"),
    @(4, "GRANULES", "Granules India Limited",      "'532482", 2.94,  517.8,   913932,    "day", "22/07/2024 11:36:32"),
    @(5, "VEDL",      "Vedanta Limited",            "'500295", 2.04,  448.75,  17818304,  "day", "22/07/2024 11:36:32")
)

$startRow = 195
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $dayWs.Cells.Item($r, 1).Value = $row[0]
    $dayWs.Cells.Item($r, 2).Value = $row[1]
    $dayWs.Cells.Item($r, 3).Value = $row[2]
    $dayWs.Cells.Item($r, 4).Value = $row[3]
    $dayWs.Cells.Item($r, 5).Value = $row[4]
    $dayWs.Cells.Item($r, 6).Value = $row[5]
    $dayWs.Cells.Item($r, 7).Value = $row[6]
    $dayWs.Cells.Item($r, 8).Value = $row[7]
    $dayWs.Cells.Item($r, 9).Value = $row[8]
}

# --- 2. "week" sheet: re-store D97:D102 (bsecode) as numbers --------------
$weekWs = $wb.Worksheets.Item("week")

$bsecodes = @{
    97  = 539448
    98  = 532343
    99  = 532478
    100 = 500483
    101 = 533273
    102 = 532432
}

foreach ($r in $bsecodes.Keys) {
    $weekWs.Cells.Item($r, 4).Value = $bsecodes[$r]
}
